$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.075.46'
$ws.Range('E2').Value = '  -3.70%  '
$ws.Range('D3').Value = '3.331.63'
$ws.Range('E3').Value = '  -5.71%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '550.67'
$ws.Range('E5').Value = '  -4.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.50'
$ws.Range('E6').Value = '  -4.38%  '
$ws.Range('E7').Value = '  -3.54%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '3.325.28'
$ws.Range('E9').Value = '  -5.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.621'
$ws.Range('E10').Value = '  -2.98%  '
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.12'
$ws.Range('E12').Value = '  -4.94%  '
$ws.Range('E13').Value = '  -2.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.99'
$ws.Range('E14').Value = '  -3.62%  '
$ws.Range('D15').Value = '3.861.65'
$ws.Range('E15').Value = '  -5.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.30'
$ws.Range('E16').Value = '  -0.95%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.117'
$ws.Range('E17').Value = '  -3.87%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.327.68'
$ws.Range('E18').Value = '  -5.61%  '
$ws.Range('E19').Value = '  -2.80%  '
$ws.Range('D20').Value = '63.979.94'
$ws.Range('E20').Value = '  -3.88%  '
$ws.Range('E21').Value = '  -3.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '423.09'
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.75'
$ws.Range('E23').Value = '  +7.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.06'
$ws.Range('E24').Value = '  -3.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.89'
$ws.Range('E25').Value = '  -2.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.17'
$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('E27').Value = '  -5.33%  '
$ws.Range('E28').Value = '  -2.80%  '
$ws.Range('E29').Value = '  -6.19%  '
$ws.Range('E30').Value = '  -3.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.67'
$ws.Range('E31').Value = '  +1.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '594.24'
$ws.Range('E32').Value = '  -8.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.37'
$ws.Range('E33').Value = '  -3.52%  '
$ws.Range('E34').Value = '  -4.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.12'
$ws.Range('E35').Value = '  -2.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.141'
$ws.Range('E37').Value = '  -10.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.52'
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.33'
$ws.Range('E39').Value = '  -6.07%  '
$ws.Range('D40').Value = '0.0₃0746'
$ws.Range('E40').Value = '  -8.95%  '
$ws.Range('E41').Value = '  -5.70%  '
$ws.Range('D42').Value = '3.090.68'
$ws.Range('E42').Value = '  -5.55%  '
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.78'
$ws.Range('E44').Value = '  -6.49%  '
$ws.Range('E45').Value = '  -4.57%  '
$ws.Range('E46').Value = '  -4.57%  '
$ws.Range('E47').Value = '  -4.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.128'
$ws.Range('E48').Value = '  -3.73%  '
$ws.Range('E49').Value = '  -4.66%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '135.45'
$ws.Range('E50').Value = '  -2.53%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.16'
$ws.Range('E51').Value = '  -6.48%  '
